$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.687.91"
$ws.Range("E2").Value = "  -1.31%  "
$ws.Range("D3").Value = "2.097.32"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "342.92"
$ws.Range("E5").Value = "  -2.04%  "
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5170"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4376"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.53"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09191"
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.166"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("E12").Value = "  -4.64%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.771"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.050.29"
$ws.Range("E14").Value = "  -2.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.152"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.56"
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.04"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06665"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.204"
$ws.Range("E22").Value = "  -1.06%  "
$ws.Range("D23").Value = "29.748.52"
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.57"
$ws.Range("E24").Value = "  -2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.303"
$ws.Range("E25").Value = "  -2.15%  "
$ws.Range("D26").Value = "2.308.69"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.88"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.86"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.496"
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.47"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.127"
$ws.Range("E31").Value = "  -4.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.690"
$ws.Range("E32").Value = "  +2.87%  "
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.190"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.956"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.353"
$ws.Range("E36").Value = "  +7.09%  "
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02578"
$ws.Range("E38").Value = "  -0.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06706"
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6983"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.334"
$ws.Range("E41").Value = "  +6.21%  "
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2217"
$ws.Range("E43").Value = "  -4.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6781"
$ws.Range("E44").Value = "  +5.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.28"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.322"
$ws.Range("E46").Value = "  -0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000360"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.619"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.204"
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.215"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.10"
$ws.Range("E51").Value = "  -2.97%  "
